# Update template email bulk excel
#
# The "account number" / "email address" columns were swapped: what used
# to be column A (account number / account numbers) now lives in column B,
# and what used to be column B (email address / email addresses) now lives
# in column A. The numeric-looking account numbers in column A were stored
# with a "quote prefix" (force-text) cell format; that format follows the
# value over to column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows (1-based) whose column-A value is text-quoted (quote-prefix format)
# in the original layout; that formatting needs to follow the value to
# column B after the swap.
$quotePrefixRows = @(2, 3, 4)

for ($r = 1; $r -le 4; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)

    $valA = $cellA.Value2
    $valB = $cellB.Value2

    $needsQuotePrefix = $quotePrefixRows -contains $r
    $prefix = if ($needsQuotePrefix) { "'" } else { "" }

    # Old column B never carried the quote-prefix format, so just drop the
    # old column A value straight into A.
    $cellA.Value2 = $valB

    # Old column A value moves into B, re-applying the quote prefix (if any)
    # so column B keeps it formatted as text, matching the original A cell.
    $cellB.Value2 = "$prefix$valA"
}

# The sheet's selection moves from the single cell B4 to the whole of
# column A's data range (A1:A4).
$ws.Range("A1:A4").Select()
